$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D text values (prices) are not auto-converted to numbers by Excel
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.464.41"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "1.854.86"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "244.62"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "0.6958"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("D7").Value = "0.9999"
$ws.Range("D8").Value = "0.07682"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.3066"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "0.07773"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "5.151"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").Value = "1.857.22"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").Value = "90.97"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "0.6917"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "6.268"
$ws.Range("E16").Value = "  -2.76%  "
$ws.Range("D17").Value = "29.455.81"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").Value = "0.000008332"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").Value = "2.099.30"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").Value = "238.15"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "0.9994"
$ws.Range("D23").Value = "7.608"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").Value = "0.9998"
$ws.Range("D25").Value = "0.1494"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").Value = "159.89"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("D27").Value = "8.884"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").Value = "18.23"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "1.533"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "1.203"
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("D33").Value = "0.05094"
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("D34").Value = "0.7734"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "1.880"
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "1.314.13"
$ws.Range("E38").Value = "  +7.78%  "
$ws.Range("D39").Value = "0.01873"
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").Value = "0.9484"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").Value = "106.15"
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("D43").Value = "5.773"
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "9.851"
$ws.Range("E45").Value = "  +2.94%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.00000000125"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "2.000.83"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.5217"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.786"
$ws.Range("E49").Value = "  +2.26%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "62.93"
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "6.949"
$ws.Range("E51").Value = "  +0.75%  "
